$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D
$ws.Range("D1").Value = "COLORE"

# Colors per row (matching the party/list in columns A/B)
# Values are entered in the same order the author typed them, so that
# the shared-string table ends up with the same ordering as the target file.
$ws.Range("D2").Value = "#FF0000"
$ws.Range("D3").Value = "#00FF00"
$ws.Range("D4").Value = "#FFA500"
$ws.Range("D6").Value = "#FFFF00"
$ws.Range("D7").Value = "#ff268f"
$ws.Range("D5").Value = "#26ffba"
$ws.Range("D8").Value = "#2945e3"
$ws.Range("D9").Value = "#00d5ff"
$ws.Range("D10").Value = "#001f9c"
$ws.Range("D11").Value = "#3086db"

# Update selection to match final state
$ws.Range("C3").Select()
